$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$nm = $p.NotesMaster
Write-Host "SlideMaster.Theme:" $sm.Theme
Write-Host "NotesMaster.Theme:" $nm.Theme
$sm.Theme | Get-Member
